$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 534.1
$ws.Range("I33").Value = 260.1111
$ws.Range("K33").Value = 260.1111
$ws.Range("M33").Value = -31.11110000000002
$ws.Range("H40").Value = 3982.2222
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 4167.5
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 4167.5
$ws.Range("M40").Value = -2325
$ws.Range("N40").Value = -4517.5
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
$ws.Range("H53").Value = 337.25
$ws.Range("I53").Value = 503.33334
$ws.Range("J53").Value = 171.16667
$ws.Range("K53").Value = 503.33334
$ws.Range("L53").Value = 171.16667
$ws.Range("M53").Value = 133.66666
$ws.Range("N53").Value = -1445.16667
$ws.Range("H61").Value = 191.8
$ws.Range("I61").Value = 191.8
$ws.Range("K61").Value = 575.4000000000001
$ws.Range("M61").Value = -403.4000000000001
$ws.Range("H88").Value = 614
$ws.Range("I88").Value = 500
$ws.Range("J88").Value = 633
$ws.Range("K88").Value = 500
$ws.Range("L88").Value = 633
$ws.Range("M88").Value = -94
$ws.Range("N88").Value = -1445
$ws.Range("H91").Value = 614
$ws.Range("I91").Value = 500
$ws.Range("J91").Value = 633
$ws.Range("K91").Value = 500
$ws.Range("L91").Value = 633
$ws.Range("M91").Value = 904
$ws.Range("N91").Value = -3441
$ws.Range("H107").Value = 1435.3429
$ws.Range("I107").Value = 623.8889
$ws.Range("J107").Value = 4174
$ws.Range("K107").Value = 623.8889
$ws.Range("L107").Value = 4174
$ws.Range("M107").Value = 1296.1111
$ws.Range("N107").Value = -8014
$ws.Range("H110").Value = 35236
$ws.Range("J110").Value = 35236
$ws.Range("L110").Value = 35236
$ws.Range("N110").Value = -43416
$ws.Range("H138").Value = 3610
$ws.Range("J138").Value = 3711.6667
$ws.Range("L138").Value = 11135.0001
$ws.Range("N138").Value = -21415.0001
$ws.Range("H141").Value = 4731.778
$ws.Range("I141").Value = 4731.778
$ws.Range("K141").Value = 14195.334
$ws.Range("M141").Value = -9015.334000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 11397.4
$ws.Range("J130").Value = 11397.4
$ws.Range("L130").Value = 11397.4
$ws.Range("N130").Value = -21437.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5833.3335
$ws.Range("I20").Value = 5833.3335
$ws.Range("K20").Value = 5833.3335
$ws.Range("M20").Value = -5586.3335
$ws.Range("H76").Value = 18036.428
$ws.Range("J76").Value = 18036.428
$ws.Range("L76").Value = 18036.428
$ws.Range("N76").Value = -18666.428
$ws.Range("H79").Value = 18036.428
$ws.Range("J79").Value = 18036.428
$ws.Range("L79").Value = 18036.428
$ws.Range("N79").Value = -20220.428
$ws.Range("H86").Value = 4500
$ws.Range("I86").Value = 4500
$ws.Range("K86").Value = 4500
$ws.Range("M86").Value = -3377
$ws.Range("H89").Value = 4500
$ws.Range("I89").Value = 4500
$ws.Range("K89").Value = 22500
$ws.Range("M89").Value = -16884
$ws.Range("H105").Value = 3034.3076
$ws.Range("I105").Value = 1868.25
$ws.Range("J105").Value = 4900
$ws.Range("K105").Value = 1868.25
$ws.Range("L105").Value = 4900
$ws.Range("M105").Value = -121.25
$ws.Range("N105").Value = -8394
$ws.Range("H134").Value = 6949.074
$ws.Range("I134").Value = 7133.12
$ws.Range("J134").Value = 4648.5
$ws.Range("K134").Value = 21399.36
$ws.Range("L134").Value = 13945.5
$ws.Range("M134").Value = -18864.36
$ws.Range("N134").Value = -19015.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 59990
$ws.Range("J109").Value = 59990
$ws.Range("L109").Value = 59990
$ws.Range("N109").Value = -62070
$ws.Range("H120").Value = 21387.166
$ws.Range("J120").Value = 22775.334
$ws.Range("L120").Value = 22775.334
$ws.Range("N120").Value = -30033.334
$ws.Range("H132").Value = 1333.3334
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 5048.5454
$ws.Range("I26").Value = 410
$ws.Range("J26").Value = 6079.3335
$ws.Range("K26").Value = 1230
$ws.Range("L26").Value = 18238.0005
$ws.Range("M26").Value = -942
$ws.Range("N26").Value = -18814.0005
$ws.Range("H86").Value = 437
$ws.Range("I86").Value = 296.25
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 888.75
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = 297.25
$ws.Range("N86").Value = -5372
$ws.Range("H89").Value = 437
$ws.Range("I89").Value = 296.25
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 2666.25
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = 3261.75
$ws.Range("N89").Value = -20856
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = ""
$ws.Range("H111").Value = 900
$ws.Range("I111").Value = 900
$ws.Range("K111").Value = 2700
$ws.Range("M111").Value = 367
$ws.Range("H138").Value = 1095
$ws.Range("I138").Value = 642.5
$ws.Range("J138").Value = 2000
$ws.Range("K138").Value = 1927.5
$ws.Range("L138").Value = 6000
$ws.Range("M138").Value = 3212.5
$ws.Range("N138").Value = -16280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6099.8335
$ws.Range("I70").Value = 5524.75
$ws.Range("J70").Value = 7250
$ws.Range("K70").Value = 5524.75
$ws.Range("L70").Value = 7250
$ws.Range("M70").Value = -5254.75
$ws.Range("N70").Value = -7790
$ws.Range("H73").Value = 6099.8335
$ws.Range("I73").Value = 5524.75
$ws.Range("J73").Value = 7250
$ws.Range("K73").Value = 5524.75
$ws.Range("L73").Value = 7250
$ws.Range("M73").Value = -4588.75
$ws.Range("N73").Value = -9122
$ws.Range("H113").Value = 955
$ws.Range("I113").Value = 955
$ws.Range("K113").Value = 955
$ws.Range("M113").Value = 1215

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3449.6667
$ws.Range("I16").Value = 3449.6667
$ws.Range("K16").Value = 3449.6667
$ws.Range("M16").Value = -3279.6667
$ws.Range("H40").Value = 5921.846
$ws.Range("I40").Value = 6107.8184
$ws.Range("K40").Value = 6107.8184
$ws.Range("M40").Value = -5971.8184
$ws.Range("H46").Value = 2604.3635
$ws.Range("I46").Value = 300
$ws.Range("K46").Value = 300
$ws.Range("M46").Value = -112
$ws.Range("H68").Value = 1000
$ws.Range("J68").Value = 1000
$ws.Range("L68").Value = 1000
$ws.Range("N68").Value = -2498
$ws.Range("H71").Value = 1000
$ws.Range("J71").Value = 1000
$ws.Range("L71").Value = 5000
$ws.Range("N71").Value = -12488
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = ""
$ws.Range("N100").Value = ""
$ws.Range("H135").Value = 48000
$ws.Range("J135").Value = 48000
$ws.Range("L135").Value = 48000
$ws.Range("N135").Value = -58140
$ws.Range("H136").Value = 3630.8572
$ws.Range("I136").Value = 3502.2
$ws.Range("K136").Value = 10506.6
$ws.Range("M136").Value = -7956.599999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 8000
$ws.Range("J82").Value = 8000
$ws.Range("L82").Value = 8000
$ws.Range("N82").Value = -8766
$ws.Range("H85").Value = 8000
$ws.Range("J85").Value = 8000
$ws.Range("L85").Value = 8000
$ws.Range("N85").Value = -10652
$ws.Range("H96").Value = 1466.3334
$ws.Range("I96").Value = 1000
$ws.Range("K96").Value = 1000
$ws.Range("M96").Value = 373
$ws.Range("H100").Value = 694.25
$ws.Range("I100").Value = 592.3333
$ws.Range("K100").Value = 1184.6666
$ws.Range("M100").Value = -643.6666
$ws.Range("H101").Value = 18367
$ws.Range("J101").Value = 18367
$ws.Range("L101").Value = 18367
$ws.Range("N101").Value = -24857
$ws.Range("H109").Value = 70000
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72774
